# Update "想去人数" (number of people interested) figures to the latest
# scraped values, matching the output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4407
$ws1.Range("F3").Value = 2465
$ws1.Range("F6").Value = 54
$ws1.Range("F7").Value = 60
$ws1.Range("F10").Value = 157
$ws1.Range("F12").Value = 1654
$ws1.Range("F14").Value = 3541

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4407
$ws4.Range("F3").Value = 2465
$ws4.Range("F7").Value = 54
$ws4.Range("F8").Value = 60
$ws4.Range("F9").Value = 44
$ws4.Range("F12").Value = 157
$ws4.Range("F16").Value = 1654
$ws4.Range("F18").Value = 3542
